$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day labels replacing the plain numbers in row 3 (A3:G3) with text
# describing progress on that day (pandas / file handling / turtle project).
# Order matches the original authoring order of the shared strings table.
$ws.Range("G3").Value = "pandas 27"
$ws.Range("F3").Value = "file handling26"
$ws.Range("A3").Value = "21 turtle"
$ws.Range("B3").Value = "22turtle"
$ws.Range("C3").Value = "23turtle"
$ws.Range("D3").Value = "24turtle"
$ws.Range("E3").Value = "25turtle"

# Match the "highlighted" style already used for similar text cells (e.g. A2:F2)
$srcStyle = $ws.Range("A2")
$target = $ws.Range("A3:G3")
$target.Interior.Color = $srcStyle.Interior.Color
$target.Borders.LineStyle = $srcStyle.Borders.LineStyle
$target.HorizontalAlignment = $srcStyle.HorizontalAlignment

# Update the active selection to J2 as in the saved workbook
$ws.Range("J2").Select()
